$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the publication date labels in row 9 (G9, H9)
$ws.Range("G9").Value = "1402-04-14 (9)"
$ws.Range("H9").Value = "1402-04-14 (2)"

# Update balance-sheet figures in column H (latest period)
$ws.Range("H14").Value = 33558574
$ws.Range("H18").Value = 64652484
$ws.Range("H27").Value = 116069676
$ws.Range("H35").Value = 56308173
$ws.Range("H37").Value = 73590840
$ws.Range("H43").Value = 77168184
$ws.Range("H56").Value = 36689608
$ws.Range("H57").Value = 38901492
$ws.Range("H58").Value = 116069676
